$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IGLD")

$ws.Range("D8").Value = 2700800
$ws.Range("E8").Value = 2782200
$ws.Range("F8").Value = 2754900
$ws.Range("G8").Value = 2498300
$ws.Range("H8").Value = 2638400
$ws.Range("I8").Value = 2835700
$ws.Range("J8").Value = 3138600

$ws.Range("D14").Value = 41900
$ws.Range("E14").Value = 26500
$ws.Range("F14").Value = 32300
$ws.Range("G14").Value = -109300
$ws.Range("H14").Value = 30900
$ws.Range("I14").Value = 11600
$ws.Range("J14").Value = 100700

$ws.Range("D15").Value = 584100
$ws.Range("E15").Value = 596200
$ws.Range("F15").Value = 587900
$ws.Range("G15").Value = 516800
$ws.Range("H15").Value = 555700
$ws.Range("I15").Value = 653100
$ws.Range("J15").Value = 823300

$ws.Range("D17").Value = 2258200
$ws.Range("E17").Value = 2268700
$ws.Range("F17").Value = 2199500
$ws.Range("G17").Value = 1787800
$ws.Range("H17").Value = 2077800
$ws.Range("I17").Value = 2301800
$ws.Range("J17").Value = 2730600

$ws.Range("D18").Value = 442500
$ws.Range("E18").Value = 513400
$ws.Range("F18").Value = 555400
$ws.Range("G18").Value = 710400
$ws.Range("H18").Value = 560600
$ws.Range("I18").Value = 533900
$ws.Range("J18").Value = 408100

$ws.Range("D20").Value = -21800
$ws.Range("E20").Value = -15500
$ws.Range("F20").Value = 30900
$ws.Range("G20").Value = 8800
$ws.Range("H20").Value = 16300
$ws.Range("I20").Value = 26200
$ws.Range("J20").Value = -16800

$ws.Range("D21").Value = 1003600
$ws.Range("E21").Value = 1093000
$ws.Range("F21").Value = 1173000
$ws.Range("G21").Value = 1234900
$ws.Range("H21").Value = 1131400
$ws.Range("I21").Value = 1211800
$ws.Range("J21").Value = 1204800

$ws.Range("D22").Value = 138800
$ws.Range("E22").Value = 254900
$ws.Range("F22").Value = 191800
$ws.Range("G22").Value = 247200
$ws.Range("H22").Value = 195100
$ws.Range("I22").Value = 208300
$ws.Range("J22").Value = 202800

$ws.Range("D23").Value = 282000
$ws.Range("E23").Value = 243100
$ws.Range("F23").Value = 394500
$ws.Range("G23").Value = 472100
$ws.Range("H23").Value = 381800
$ws.Range("I23").Value = 351800
$ws.Range("J23").Value = 188400

$ws.Range("D24").Value = 95700
$ws.Range("E24").Value = 121900
$ws.Range("F24").Value = 95700
$ws.Range("G24").Value = 184000
$ws.Range("H24").Value = 144600
$ws.Range("I24").Value = 153400
$ws.Range("J24").Value = 180200

$ws.Range("D26").Value = 186200
$ws.Range("E26").Value = 121100
$ws.Range("F26").Value = 298800
$ws.Range("G26").Value = 288000
$ws.Range("H26").Value = 237300
$ws.Range("I26").Value = 198400
$ws.Range("J26").Value = 8300

$ws.Range("D27").Value = -4100
$ws.Range("E27").Value = -55700
$ws.Range("F27").Value = 24000
$ws.Range("G27").Value = -28400
$ws.Range("H27").Value = 7200
$ws.Range("I27").Value = -10200
$ws.Range("J27").Value = -73400

$ws.Range("D32").Value = 21800
$ws.Range("E32").Value = 15500
$ws.Range("F32").Value = -30900
$ws.Range("G32").Value = -8800
$ws.Range("H32").Value = -16300
$ws.Range("I32").Value = -26200
$ws.Range("J32").Value = 16800

$ws.Range("D33").Value = -4100
$ws.Range("E33").Value = -55700
$ws.Range("F33").Value = 24000
$ws.Range("G33").Value = -28400
$ws.Range("H33").Value = 7200
$ws.Range("I33").Value = -10200
$ws.Range("J33").Value = -73400

$ws.Range("D35").Value = -4100
$ws.Range("E35").Value = -55700
$ws.Range("F35").Value = 24000
$ws.Range("G35").Value = -28400
$ws.Range("H35").Value = 7200
$ws.Range("I35").Value = -10200
$ws.Range("J35").Value = -73400

$ws.Range("D41").Value = 664400
$ws.Range("E41").Value = 223500
$ws.Range("F41").Value = 170800
$ws.Range("G41").Value = 202000
$ws.Range("H41").Value = 239200
$ws.Range("I41").Value = 210800
$ws.Range("J41").Value = 798500

$ws.Range("D42").Value = 212200
$ws.Range("E42").Value = 342100
$ws.Range("F42").Value = 489400
$ws.Range("G42").Value = 939700
$ws.Range("H42").Value = 515400
$ws.Range("I42").Value = 456600
$ws.Range("J42").Value = 427100

$ws.Range("D43").Value = 596500
$ws.Range("E43").Value = 571700
$ws.Range("F43").Value = 605000
$ws.Range("G43").Value = 652800
$ws.Range("H43").Value = 887600
$ws.Range("I43").Value = 1766300
$ws.Range("J43").Value = 1830900

$ws.Range("D44").Value = 34500
$ws.Range("E44").Value = 29200
$ws.Range("F44").Value = 31700
$ws.Range("G44").Value = 26500
$ws.Range("H44").Value = 32300
$ws.Range("I44").Value = 33900
$ws.Range("J44").Value = 56300

$ws.Range("D45").Value = 18200
$ws.Range("E45").Value = 40000
$ws.Range("F45").Value = 84400
$ws.Range("G45").Value = 60700
$ws.Range("H45").Value = 95200
$ws.Range("I45").Value = 75600
$ws.Range("J45").Value = 72300

$ws.Range("D46").Value = 1525700
$ws.Range("E46").Value = 1206500
$ws.Range("F46").Value = 1377000
$ws.Range("G46").Value = 1881600
$ws.Range("H46").Value = 1673900
$ws.Range("I46").Value = 1644900
$ws.Range("J46").Value = 1853800

$ws.Range("D47").Value = 139100
$ws.Range("E47").Value = 182600
$ws.Range("F47").Value = 192900
$ws.Range("G47").Value = 522600
$ws.Range("H47").Value = 482300
$ws.Range("I47").Value = 598400
$ws.Range("J47").Value = 738600

$ws.Range("D48").Value = 1914700
$ws.Range("E48").Value = 1951200
$ws.Range("F48").Value = 1990100
$ws.Range("G48").Value = 1813200
$ws.Range("H48").Value = 3609300
$ws.Range("I48").Value = 2742200
$ws.Range("J48").Value = 4003600

$ws.Range("D49").Value = 1736500
$ws.Range("E49").Value = 1921900
$ws.Range("F49").Value = 2622200
$ws.Range("G49").Value = 1630000
$ws.Range("H49").Value = 3649100
$ws.Range("I49").Value = 4001700
$ws.Range("J49").Value = 6692000

$ws.Range("D52").Value = 432100
$ws.Range("E52").Value = 401200
$ws.Range("F52").Value = 533300
$ws.Range("G52").Value = 100400
$ws.Range("H52").Value = 121700
$ws.Range("I52").Value = 141300
$ws.Range("J52").Value = 175200

$ws.Range("D54").Value = 5748100
$ws.Range("E54").Value = 5663400
$ws.Range("F54").Value = 6182900
$ws.Range("G54").Value = 5947900
$ws.Range("H54").Value = 5907000
$ws.Range("I54").Value = 6292200
$ws.Range("J54").Value = 6969000

$ws.Range("D57").Value = 478700
$ws.Range("E57").Value = 458300
$ws.Range("F57").Value = 473700
$ws.Range("G57").Value = 183200
$ws.Range("H57").Value = 198900
$ws.Range("I57").Value = 218800
$ws.Range("J57").Value = 246100

$ws.Range("D58").Value = 539400
$ws.Range("E58").Value = 601700
$ws.Range("F58").Value = 612200
$ws.Range("G58").Value = 430700
$ws.Range("H58").Value = 432100
$ws.Range("I58").Value = 471000
$ws.Range("J58").Value = 364200

$ws.Range("D59").Value = 147300
$ws.Range("E59").Value = 155900
$ws.Range("F59").Value = 390700
$ws.Range("G59").Value = 482500
$ws.Range("H59").Value = 501300
$ws.Range("I59").Value = 661900
$ws.Range("J59").Value = 938900

$ws.Range("D60").Value = 1165400
$ws.Range("E60").Value = 1215900
$ws.Range("F60").Value = 1476600
$ws.Range("G60").Value = 1096400
$ws.Range("H60").Value = 1132300
$ws.Range("I60").Value = 1351600
$ws.Range("J60").Value = 1320700

$ws.Range("D61").Value = 3627800
$ws.Range("E61").Value = 3377300
$ws.Range("F61").Value = 3646000
$ws.Range("G61").Value = 3702300
$ws.Range("H61").Value = 3510800
$ws.Range("I61").Value = 3552200
$ws.Range("J61").Value = 3800500

$ws.Range("D62").Value = 277300
$ws.Range("E62").Value = 315100
$ws.Range("F62").Value = 342700
$ws.Range("G62").Value = 386000
$ws.Range("H62").Value = 394000
$ws.Range("I62").Value = 428200
$ws.Range("J62").Value = 677600

$ws.Range("D66").Value = 5699300
$ws.Range("E66").Value = 5609900
$ws.Range("F66").Value = 6208600
$ws.Range("G66").Value = 5998300
$ws.Range("H66").Value = 5930700
$ws.Range("I66").Value = 6317600
$ws.Range("J66").Value = 6976400

$ws.Range("D72").Value = -86100
$ws.Range("E72").Value = -81400
$ws.Range("F72").Value = -160600
$ws.Range("G72").Value = -185400
$ws.Range("H72").Value = -158600
$ws.Range("I72").Value = -160300
$ws.Range("J72").Value = -142400

$ws.Range("D76").Value = 48800
$ws.Range("E76").Value = 53500
$ws.Range("F76").Value = -25700
$ws.Range("G76").Value = -50500
$ws.Range("H76").Value = -23700
$ws.Range("I76").Value = -25400
$ws.Range("J76").Value = -7400

$ws.Range("D81").Value = -4100
$ws.Range("E81").Value = -55700
$ws.Range("F81").Value = 24000
$ws.Range("G81").Value = -28400
$ws.Range("H81").Value = 7200
$ws.Range("I81").Value = -10200
$ws.Range("J81").Value = -73400

$ws.Range("D83").Value = 584100
$ws.Range("E83").Value = 596200
$ws.Range("F83").Value = 587900
$ws.Range("G83").Value = 516800
$ws.Range("H83").Value = 555700
$ws.Range("I83").Value = 653100
$ws.Range("J83").Value = 815300

$ws.Range("D89").Value = 960100
$ws.Range("E89").Value = 953800
$ws.Range("F89").Value = 1005400
$ws.Range("G89").Value = 1043500
$ws.Range("H89").Value = 1142500
$ws.Range("I89").Value = 1102800
$ws.Range("J89").Value = 875200

$ws.Range("D91").Value = -312000
$ws.Range("E91").Value = -329100
$ws.Range("F91").Value = -365300
$ws.Range("G91").Value = -298200
$ws.Range("H91").Value = -287500
$ws.Range("I91").Value = -350700
$ws.Range("J91").Value = -427100

$ws.Range("D94").Value = -266500
$ws.Range("E94").Value = -285800
$ws.Range("F94").Value = 105100
$ws.Range("G94").Value = -548800
$ws.Range("H94").Value = -285000
$ws.Range("I94").Value = -320000
$ws.Range("J94").Value = -572500

$ws.Range("D100").Value = -252700
$ws.Range("E100").Value = -615300
$ws.Range("F100").Value = -1141700
$ws.Range("G100").Value = -531900
$ws.Range("H100").Value = -829100
$ws.Range("I100").Value = -971200
$ws.Range("J100").Value = -14900

$ws.Range("D102").Value = 440900
$ws.Range("E102").Value = 52700
$ws.Range("F102").Value = -31200
$ws.Range("G102").Value = -37200
$ws.Range("H102").Value = 28400
$ws.Range("I102").Value = -188400
$ws.Range("J102").Value = 287800
